# Rename the "individualParamsFile" row to "individualsFile" / "Individuals.xlsx"
# (row 5), and remove the now-obsolete "individualPhysiologyFile" /
# "IndividualBiometrics.xlsx" row (row 6) - individual physiology info is now
# part of the "Individuals" file. All rows below shift up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "individualsFile"
$ws.Range("B5").Value = "Individuals.xlsx"

$ws.Rows("6:6").Delete()
